$wb = $excel.ActiveWorkbook

# "Repayment Schedule" sheet gets a new blank spacer column inserted before
# the old column N ("Late" header), pushing Late/Outstanding one column right
# and mirroring the existing spacer pattern used for column E / old column O.
$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")
$wsSchedule.Columns("N").Insert()

# Move the active tab from "Transactions" to "Repayment Schedule", and move
# the in-sheet selection on "Repayment Schedule" to K8.
$wsSchedule.Activate()
$wsSchedule.Range("K8").Select()
